$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($null -ne $val) {
        $scaled = [Math]::Round($val * 10000, 4)
        $cell.Value2 = $scaled
    }
}
